$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped from coinranking.com, refreshed by the sync job.
# Column E (% change) values are always text (percent sign, padded with
# spaces) so a plain .Value assignment keeps them as text.
# Column D (price) values are sometimes plain decimal numbers; assigning
# those via .Value would make Excel auto-coerce them to the Number type
# (dropping trailing zeros / using scientific notation), so for those we
# force the Text number format first, then restore the default "Normal"
# cell style once the literal text is stored.

$ws.Range("D2").Value = '56.315.19'
$ws.Range("E2").Value = '  -4.08%  '
$ws.Range("D3").Value = '2.375.48'
$ws.Range("E3").Value = '  -4.46%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.61%  '
$ws.Range("D9").Value = '2.397.72'
$ws.Range("E9").Value = '  -3.34%  '
$ws.Range("E10").Value = '  -3.96%  '
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("E12").Value = '  -9.02%  '
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").Value = '2.803.56'
$ws.Range("E14").Value = '  -4.49%  '
$ws.Range("D15").Value = '56.228.30'
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.44%  '
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("D18").Value = '2.388.55'
$ws.Range("E18").Value = '  -4.70%  '
$ws.Range("E19").Value = '  -4.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '308.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.63%  '
$ws.Range("E21").Value = '  -4.75%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '2.488.38'
$ws.Range("E27").Value = '  -7.59%  '
$ws.Range("E28").Value = '  -5.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("E31").Value = '  -5.09%  '
$ws.Range("E32").Value = '  -3.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("E35").Value = '  -6.71%  '
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  -4.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.783'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("E44").Value = '  -3.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '251.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.563'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0897'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("E49").Value = '  -4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("E51").Value = '  -4.74%  '
